$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 00:44:05"
$wsZhCn.Range("H2").Value = "2016-03-22 00:44:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 00:44:08"
$wsDeDe.Range("H2").Value = "2016-03-22 00:44:39"
